$wb = $excel.ActiveWorkbook

# Excel's ColumnWidth (character units) is stored in the XML "width" attribute
# with a constant +5/6 (0.8333333333333336) padding offset baked in by this
# runtime, so to land on a stored width of exactly 40 we request 39.1666...
$targetColWidth = 40 - 0.8333333333333336

# 1. Overview sheet: the Status text for the 3a0fe57b file ("Ready for
#    handoff") becomes "Handback transform failed". The very same text is
#    also shown in the zh-cn / de-de per-language tables' Status column for
#    that file, so update those too.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Handback transform failed"
$overview.Range("F3").Value = "Handback transform failed"

# 2. zh-cn sheet: Status + new Error Detail message for the 3a0fe57b row,
#    plus widen the Error Detail column.
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Handback transform failed"
$zhcn.Range("P3").Value = "Handback file name: ttr2zh21.t0n is different with handoff file name: 3a0fe57b-2c82-4ef9-a76e-b7849e5353a9.f8e141b8b7ba2a27138d84619df911961e74b0eb.zh-cn."
$zhcn.Columns.Item(16).ColumnWidth = $targetColWidth

# 3. de-de sheet: same kind of update, with the de-de handoff file name.
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Handback transform failed"
$dede.Range("P3").Value = "Handback file name: ttr2zh21.t0n is different with handoff file name: 3a0fe57b-2c82-4ef9-a76e-b7849e5353a9.f8e141b8b7ba2a27138d84619df911961e74b0eb.de-de."
$dede.Columns.Item(16).ColumnWidth = $targetColWidth
